$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 510000017
$ws.Range("B1").Value = 1212
$ws.Range("A2:B2").EntireRow.Delete()
$ws.Outline.ShowLevels(0, 1)
$ws.Range("B2").Select()
